# Regenerate s_val data to filter save games: update numeric values in
# columns B-E and the summed G column for data rows 2-7 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; B = 3.230985683306322;  C = 0.3127903958511391; D = 3.900430680208489;  E = 0.496779210170732;  G = 7.940985969536682 },
    @{ Row = 3; B = 0.3048080303191223; C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 2.626907116734944 },
    @{ Row = 4; B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 },
    @{ Row = 5; B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 },
    @{ Row = 6; B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 8.660232485948974;  G = 17.45944343273191 },
    @{ Row = 7; B = 3.230985683306322;  C = 10.29869402782916;  D = 3.900430680208489;  E = 8.660232485948974;  G = 26.09034287729295 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("B$r").Value = $u.B
    $ws.Range("C$r").Value = $u.C
    $ws.Range("D$r").Value = $u.D
    $ws.Range("E$r").Value = $u.E
    $ws.Range("G$r").Value = $u.G
}
